$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("mock1")
$ws1.Range("C2").Value = 8858
$ws1.Range("D2").Value = 8859
$ws1.Range("G2").Value = 1728

$ws1.Range("C3").Value = 2281
$ws1.Range("D3").Value = 2281

$ws1.Range("C4").Value = 1514
$ws1.Range("D4").Value = 1514

$ws1.Range("C5").Value = 1232
$ws1.Range("D5").Value = 1233
$ws1.Range("G5").Value = 343

$ws1.Range("G6").Value = 259

$ws1.Range("C7").Value = 807
$ws1.Range("D7").Value = 807
$ws1.Range("G7").Value = 223

$ws1.Range("C8").Value = 1125
$ws1.Range("D8").Value = 1125
$ws1.Range("G8").Value = 160

$ws1.Range("G9").Value = 200

$ws1.Range("C10").Value = 157
$ws1.Range("D10").Value = 157

$ws1.Range("G15").Value = 4
$ws1.Range("G16").Value = 8

$ws2 = $wb.Worksheets.Item("mock2")
$ws2.Range("C2").Value = 4633
$ws2.Range("D2").Value = 4625
$ws2.Range("E2").Value = 876

$ws2.Range("C3").Value = 1960
$ws2.Range("D3").Value = 1960
$ws2.Range("E3").Value = 591

$ws2.Range("C4").Value = 1492
$ws2.Range("D4").Value = 1492

$ws2.Range("E5").Value = 160

$ws2.Range("C6").Value = 163
$ws2.Range("D6").Value = 163
$ws2.Range("E6").Value = 43

$ws2.Range("E7").Value = 52

$ws2.Range("E8").Value = 16

$ws2.Range("E11").Value = 3
